$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q1" immediately before the "总计" sheet.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(6)
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# ---------------------------------------------------------------------------
# 2) Header row (matches the style used by the other quarterly sheets:
#    bold font, thin border, centered / top-aligned).
# ---------------------------------------------------------------------------
$headerRange = $q1.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------------
# 3) Data rows 2-8 (fund holdings). Columns B, C, D, E, F, G are stored as
#    text (fund codes have leading zeros, % / money figures kept verbatim),
#    column H ("仓位排名") is numeric, column A is the same styled index
#    column used throughout the workbook.
# ---------------------------------------------------------------------------
$rows = @(
    @{ idx = 0; code = "050009"; name = "博时新兴成长混合";                     scale = "32.53"; pos = "86.97"; pct = "6.89"; mv = "2.2413"; rank = 5 },
    @{ idx = 1; code = "012428"; name = "华夏核心制造混合型证券投资基金A";      scale = "43.75"; pos = "91.55"; pct = "4.87"; mv = "2.1306"; rank = 6 },
    @{ idx = 2; code = "012429"; name = "华夏核心制造混合型证券投资基金C";      scale = "9.10";  pos = "91.55"; pct = "4.87"; mv = "0.4432"; rank = 6 },
    @{ idx = 3; code = "398011"; name = "中海分红增利混合";                     scale = "2.65";  pos = "91.97"; pct = "4.21"; mv = "0.1116"; rank = 6 },
    @{ idx = 4; code = "002213"; name = "中海顺鑫灵活配置混合";                 scale = "0.81";  pos = "91.48"; pct = "4.19"; mv = "0.0339"; rank = 6 },
    @{ idx = 5; code = "011377"; name = "创金合信积极成长股票A";                scale = "0.29";  pos = "94.90"; pct = "4.47"; mv = "0.0130"; rank = 8 },
    @{ idx = 6; code = "011378"; name = "创金合信积极成长股票C";                scale = "0.11";  pos = "94.90"; pct = "4.47"; mv = "0.0049"; rank = 8 }
)

# Force text storage for the code / numeric-looking text columns so that
# leading zeros and formatting survive, then release the format again so no
# extra style stays attached to the cells (matches the plain, un-styled data
# cells used elsewhere in the workbook).
$textCols = $q1.Range("B2:G8")
$textCols.NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    $idxCell = $q1.Range("A$r")
    $idxCell.Value = $row.idx
    $idxCell.Font.Bold = $true
    $idxCell.Borders.LineStyle = 1
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160

    $q1.Range("B$r").Value = $row.code
    $q1.Range("C$r").Value = $row.name
    $q1.Range("D$r").Value = $row.scale
    $q1.Range("E$r").Value = $row.pos
    $q1.Range("F$r").Value = $row.pct
    $q1.Range("G$r").Value = $row.mv
    $q1.Range("H$r").Value = $row.rank

    $r = $r + 1
}

# Drop the helper "@" number format back to Normal so the text cells end up
# completely unstyled, just like the rest of the data cells in this sheet.
$textCols.Style = "Normal"

# ---------------------------------------------------------------------------
# 4) Prepend a "2022-Q1" summary row to the "总计" sheet (now the 7th sheet),
#    pushing the existing rows down by one. NOTE: re-fetch the sheet by name
#    rather than reusing $totalSheet - inserting a new sheet at its old
#    position rebinds that handle to whatever now sits at that index.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$a2 = $totalSheet.Range("A2")
$a2.Value = 0
$a2.Font.Bold = $true
$a2.Borders.LineStyle = 1
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160

$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 4.98

# The row-insert only shifts position; the running index in column A keeps
# its old numeric value and needs to be bumped by one for every pre-existing
# row (they are now one row further down the list).
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
